$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.430.15"
$ws.Range("E2").Value = "  +0.07%  "

$ws.Range("D3").Value = "1.816.91"
$ws.Range("E3").Value = "  -0.46%  "

$ws.Range("E4").Value = "  +0.15%  "

$ws.Range("D5").Value = "'315.25"
$ws.Range("E5").Value = "  -0.68%  "

$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("D7").Value = "'0.5087"
$ws.Range("E7").Value = "  -4.55%  "

$ws.Range("E8").Value = "  -1.94%  "

$ws.Range("D9").Value = "'0.08105"
$ws.Range("E9").Value = "  +6.46%  "

$ws.Range("D10").Value = "'41.64"
$ws.Range("E10").Value = "  -0.45%  "

$ws.Range("D11").Value = "'1.105"
$ws.Range("E11").Value = "  -0.38%  "

$ws.Range("B12").Value = "Solana"
$ws.Range("C12").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D12").Value = "'20.98"
$ws.Range("E12").Value = "  +0.36%  "

$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'6.285"
$ws.Range("E13").Value = "  -0.82%  "

$ws.Range("E14").Value = "  +0.11%  "

$ws.Range("D15").Value = "'7.499"
$ws.Range("E15").Value = "  -1.60%  "

$ws.Range("D16").Value = "1.816.24"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").Value = "'0.00001136"
$ws.Range("E17").Value = "  +5.53%  "

$ws.Range("D18").Value = "'92.48"
$ws.Range("E18").Value = "  +3.30%  "

$ws.Range("D19").Value = "'0.06627"
$ws.Range("E19").Value = "  +0.49%  "

$ws.Range("D20").Value = "'17.69"
$ws.Range("E20").Value = "  +0.04%  "

$ws.Range("D22").Value = "'6.092"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").Value = "28.456.30"
$ws.Range("E23").Value = "  +0.15%  "

$ws.Range("E24").Value = "  +0.90%  "

$ws.Range("D25").Value = "'2.267"
$ws.Range("E25").Value = "  +2.28%  "

$ws.Range("D26").Value = "'21.16"
$ws.Range("E26").Value = "  +2.56%  "

$ws.Range("D27").Value = "'155.67"
$ws.Range("E27").Value = "  -1.27%  "

$ws.Range("D28").Value = "2.027.28"
$ws.Range("E28").Value = "  -0.60%  "

$ws.Range("D29").Value = "'2.399"
$ws.Range("E29").Value = "  -2.29%  "

$ws.Range("D30").Value = "'126.11"
$ws.Range("E30").Value = "  +1.70%  "

$ws.Range("E31").Value = "  -0.40%  "

$ws.Range("E32").Value = "  -1.43%  "

$ws.Range("D33").Value = "'5.786"
$ws.Range("E33").Value = "  +2.19%  "

$ws.Range("E34").Value = "  +0.17%  "

$ws.Range("D35").Value = "'0.07022"
$ws.Range("E35").Value = "  -5.99%  "

$ws.Range("D36").Value = "'0.2221"
$ws.Range("E36").Value = "  -0.44%  "

$ws.Range("D37").Value = "'5.217"
$ws.Range("E37").Value = "  +0.36%  "

$ws.Range("D38").Value = "'0.02329"
$ws.Range("E38").Value = "  -0.52%  "

$ws.Range("D39").Value = "'8.828"
$ws.Range("E39").Value = "  -0.92%  "

$ws.Range("E40").Value = "  +0.44%  "

$ws.Range("D41").Value = "'11.27"
$ws.Range("E41").Value = "  -0.40%  "

$ws.Range("D42").Value = "'1.173"
$ws.Range("E42").Value = "  -0.36%  "

$ws.Range("D43").Value = "'1.001"
$ws.Range("E43").Value = "  +0.13%  "

$ws.Range("D44").Value = "'1.402"
$ws.Range("E44").Value = "  +0.62%  "

$ws.Range("D45").Value = "'13.39"

$ws.Range("D46").Value = "'3.741"
$ws.Range("E46").Value = "  +1.00%  "

$ws.Range("D47").Value = "'0.5914"
$ws.Range("E47").Value = "  +1.23%  "

$ws.Range("D48").Value = "'124.86"
$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").Value = "'1.973"
$ws.Range("E49").Value = "  -1.00%  "

$ws.Range("E50").Value = "  -1.25%  "

$ws.Range("D51").Value = "'0.06886"
$ws.Range("E51").Value = "  -0.07%  "
